# Edit: "Update countries & provincias Spain"
# Refreshes the COVID-19 country snapshot on sheet "Pais":
#   - bumps the "last updated" timestamp in A1
#   - several countries swapped table positions as their totals overtook
#     their neighbours, so the country name in column A is corrected for
#     every row whose rank changed
#   - refreshes the daily case/recovery/death counters (columns B-H) for
#     every row whose figures moved

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- "Datos actualizados" timestamp (A1) ---
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 7 de Julio de 2020 a las 02:06"

# Row 4: stats updated
$ws.Cells.Item(4, 2).Value = 3039670  # Casos totales
$ws.Cells.Item(4, 3).Value = 49423  # Nuevos casos
$ws.Cells.Item(4, 4).Value = 1310649  # Casos activos
$ws.Cells.Item(4, 5).Value = 1596069  # Recuperados
$ws.Cells.Item(4, 7).Value = 351  # Muertes hoy
$ws.Cells.Item(4, 8).Value = 132952  # Muertes

# Row 5: stats updated
$ws.Cells.Item(5, 2).Value = 1626071  # Casos totales
$ws.Cells.Item(5, 3).Value = 21486  # Nuevos casos
$ws.Cells.Item(5, 5).Value = 581900  # Recuperados
$ws.Cells.Item(5, 7).Value = 656  # Muertes hoy
$ws.Cells.Item(5, 8).Value = 65556  # Muertes

# Row 23: stats updated
$ws.Cells.Item(23, 2).Value = 105934  # Casos totales
$ws.Cells.Item(23, 3).Value = 398  # Nuevos casos
$ws.Cells.Item(23, 4).Value = 69570  # Casos activos
$ws.Cells.Item(23, 5).Value = 27671  # Recuperados
$ws.Cells.Item(23, 7).Value = 9  # Muertes hoy
$ws.Cells.Item(23, 8).Value = 8693  # Muertes

# Row 26: stats updated
$ws.Cells.Item(26, 2).Value = 80447  # Casos totales
$ws.Cells.Item(26, 3).Value = 2632  # Nuevos casos
$ws.Cells.Item(26, 5).Value = 50334  # Recuperados
$ws.Cells.Item(26, 7).Value = 75  # Muertes hoy
$ws.Cells.Item(26, 8).Value = 1582  # Muertes

# Row 43: Bolivia -> Panama; stats updated
$ws.Cells.Item(43, 1).Value = "Panama"
$ws.Cells.Item(43, 2).Value = 39334  # Casos totales
$ws.Cells.Item(43, 3).Value = 1185  # Nuevos casos
$ws.Cells.Item(43, 4).Value = 18036  # Casos activos
$ws.Cells.Item(43, 5).Value = 20528  # Recuperados
$ws.Cells.Item(43, 7).Value = 23  # Muertes hoy
$ws.Cells.Item(43, 8).Value = 770  # Muertes

# Row 44: Panama -> Bolivia; stats updated
$ws.Cells.Item(44, 1).Value = "Bolivia"
$ws.Cells.Item(44, 2).Value = 39297  # Casos totales
$ws.Cells.Item(44, 3).Value = 1226  # Nuevos casos
$ws.Cells.Item(44, 4).Value = 11667  # Casos activos
$ws.Cells.Item(44, 5).Value = 26196  # Recuperados
$ws.Cells.Item(44, 7).Value = 56  # Muertes hoy
$ws.Cells.Item(44, 8).Value = 1434  # Muertes

# Row 55: Honduras -> Guatemala; stats updated
$ws.Cells.Item(55, 1).Value = "Guatemala"
$ws.Cells.Item(55, 2).Value = 23972  # Casos totales
$ws.Cells.Item(55, 3).Value = 724  # Nuevos casos
$ws.Cells.Item(55, 4).Value = 3429  # Casos activos
$ws.Cells.Item(55, 5).Value = 19562  # Recuperados
$ws.Cells.Item(55, 7).Value = 34  # Muertes hoy
$ws.Cells.Item(55, 8).Value = 981  # Muertes

# Row 56: Guatemala -> Honduras; stats updated
$ws.Cells.Item(56, 1).Value = "Honduras"
$ws.Cells.Item(56, 2).Value = 23943  # Casos totales
$ws.Cells.Item(56, 3).Value = 1022  # Nuevos casos
$ws.Cells.Item(56, 4).Value = 2490  # Casos activos
$ws.Cells.Item(56, 5).Value = 20814  # Recuperados
$ws.Cells.Item(56, 7).Value = 10  # Muertes hoy
$ws.Cells.Item(56, 8).Value = 639  # Muertes

# Row 57: Azerbaiyan -> Ghana; stats updated
$ws.Cells.Item(57, 1).Value = "Ghana"
$ws.Cells.Item(57, 2).Value = 21077  # Casos totales
$ws.Cells.Item(57, 3).Value = 992  # Nuevos casos
$ws.Cells.Item(57, 4).Value = 16070  # Casos activos
$ws.Cells.Item(57, 5).Value = 4878  # Recuperados
$ws.Cells.Item(57, 7).Value = 7  # Muertes hoy
$ws.Cells.Item(57, 8).Value = 129  # Muertes

# Row 58: Ghana -> Azerbaiyan; stats updated
$ws.Cells.Item(58, 1).Value = "Azerbaiyan"
$ws.Cells.Item(58, 2).Value = 20837  # Casos totales
$ws.Cells.Item(58, 3).Value = 513  # Nuevos casos
$ws.Cells.Item(58, 4).Value = 12182  # Casos activos
$ws.Cells.Item(58, 5).Value = 8397  # Recuperados
$ws.Cells.Item(58, 7).Value = 8  # Muertes hoy
$ws.Cells.Item(58, 8).Value = 258  # Muertes

# Row 118: stats updated
$ws.Cells.Item(118, 2).Value = 1790  # Casos totales
$ws.Cells.Item(118, 3).Value = 25  # Nuevos casos
$ws.Cells.Item(118, 4).Value = 760  # Casos activos
$ws.Cells.Item(118, 5).Value = 1005  # Recuperados

# Row 133: Ruanda -> Libia; stats updated
$ws.Cells.Item(133, 1).Value = "Libia"
$ws.Cells.Item(133, 2).Value = 1117  # Casos totales
$ws.Cells.Item(133, 3).Value = 71  # Nuevos casos
$ws.Cells.Item(133, 4).Value = 269  # Casos activos
$ws.Cells.Item(133, 5).Value = 814  # Recuperados
$ws.Cells.Item(133, 7).Value = 2  # Muertes hoy
$ws.Cells.Item(133, 8).Value = 34  # Muertes

# Row 134: Niger -> Ruanda; stats updated
$ws.Cells.Item(134, 1).Value = "Ruanda"
$ws.Cells.Item(134, 2).Value = 1113  # Casos totales
$ws.Cells.Item(134, 3).Value = 8  # Nuevos casos
$ws.Cells.Item(134, 4).Value = 575  # Casos activos
$ws.Cells.Item(134, 5).Value = 535  # Recuperados
$ws.Cells.Item(134, 8).Value = 3  # Muertes

# Row 135: Libia -> Niger; stats updated
$ws.Cells.Item(135, 1).Value = "Niger"
$ws.Cells.Item(135, 2).Value = 1093  # Casos totales
$ws.Cells.Item(135, 3).Value = 5  # Nuevos casos
$ws.Cells.Item(135, 4).Value = 968  # Casos activos
$ws.Cells.Item(135, 5).Value = 57  # Recuperados
$ws.Cells.Item(135, 8).Value = 68  # Muertes

# Row 140: stats updated
$ws.Cells.Item(140, 2).Value = 960  # Casos totales
$ws.Cells.Item(140, 3).Value = 4  # Nuevos casos
$ws.Cells.Item(140, 4).Value = 858  # Casos activos
$ws.Cells.Item(140, 5).Value = 73  # Recuperados
$ws.Cells.Item(140, 7).Value = 1  # Muertes hoy
$ws.Cells.Item(140, 8).Value = 29  # Muertes

# Row 147: Jamaica -> Zimbabue; stats updated
$ws.Cells.Item(147, 1).Value = "Zimbabue"
$ws.Cells.Item(147, 2).Value = 734  # Casos totales
$ws.Cells.Item(147, 3).Value = 18  # Nuevos casos
$ws.Cells.Item(147, 4).Value = 197  # Casos activos
$ws.Cells.Item(147, 5).Value = 528  # Recuperados
$ws.Cells.Item(147, 7).Value = 1  # Muertes hoy
$ws.Cells.Item(147, 8).Value = 9  # Muertes

# Row 148: Santo Tome y Principe -> Jamaica; stats updated
$ws.Cells.Item(148, 1).Value = "Jamaica"
$ws.Cells.Item(148, 2).Value = 732  # Casos totales
$ws.Cells.Item(148, 3).Value = 4  # Nuevos casos
$ws.Cells.Item(148, 4).Value = 584  # Casos activos
$ws.Cells.Item(148, 5).Value = 138  # Recuperados
$ws.Cells.Item(148, 8).Value = 10  # Muertes

# Row 149: Zimbabue -> Santo Tome y Principe; stats updated
$ws.Cells.Item(149, 1).Value = "Santo Tome y Principe"
$ws.Cells.Item(149, 2).Value = 721  # Casos totales
$ws.Cells.Item(149, 3).Value = 1  # Nuevos casos
$ws.Cells.Item(149, 4).Value = 269  # Casos activos
$ws.Cells.Item(149, 5).Value = 439  # Recuperados
$ws.Cells.Item(149, 8).Value = 13  # Muertes

# Row 152: stats updated
$ws.Cells.Item(152, 4).Value = 450  # Casos activos
$ws.Cells.Item(152, 5).Value = 215  # Recuperados

# Row 154: stats updated
$ws.Cells.Item(154, 2).Value = 614  # Casos totales
$ws.Cells.Item(154, 3).Value = 20  # Nuevos casos
$ws.Cells.Item(154, 5).Value = 285  # Recuperados

# Row 179: stats updated
$ws.Cells.Item(179, 4).Value = 117  # Casos activos
$ws.Cells.Item(179, 5).Value = 8  # Recuperados

# Row 209: Islas Malvinas -> Groenlandia
$ws.Cells.Item(209, 1).Value = "Groenlandia"

# Row 210: Groenlandia -> Islas Malvinas
$ws.Cells.Item(210, 1).Value = "Islas Malvinas"
